{"js": "// Change the \"Availability\" paragraph so it reads:\n//   \"Available on the web at <benjarman.uk link> or as a PDF from <repository link>.\"\n// instead of:\n//   \"Embargoed but can be requested from <repository link>.\"\n\nconst body = context.document.body;\n\n// 1) Swap the introductory text.\nconst intro = body.search(\"Embargoed but can be requested from\", { matchCase: true });\nintro.load(\"text\");\nawait context.sync();\n\nif (intro.items.length > 0) {\n  intro.items[0].insertText(\"Available on the web at\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Locate the existing (repository) hyperlink target text so we can insert the\n//    new benjarman.uk link, plus the \"or as a PDF from\" wording, right before it.\nconst repoUrl = \"https://www.repository.cam.ac.uk/handle/1810/369145\";\nconst repoResults = body.search(repoUrl, { matchCase: true });\nrepoResults.load(\"text\");\nawait context.sync();\n\nconst repoRange = repoResults.items[0];\n\n// Insert the new link text immediately before the repository link text and turn\n// it into a hyperlink pointing at the new location.\nconst newLinkRange = repoRange.insertText(\"https://benjarman.uk/phd_thesis\", Word.InsertLocation.before);\nawait context.sync();\nnewLinkRange.hyperlink = \"https://benjarman.uk/phd_thesis\";\nawait context.sync();\n\n// Insert the connecting words between the two links, still before the repository\n// link's own text. The insertion point sits right next to the hyperlink we just\n// created, so explicitly clear any inherited character style (otherwise the new\n// plain-text run would inherit the \"Hyperlink\" character style).\nconst connector = repoRange.insertText(\" or as a PDF from \", Word.InsertLocation.before);\nconnector.style = \"Default Paragraph Font\";\nawait context.sync();\n", "ps1": "# Change the \"Availability\" paragraph so it reads:\n#   \"Available on the web at <benjarman.uk link> or as a PDF from <repository link>.\"\n# instead of:\n#   \"Embargoed but can be requested from <repository link>.\"\n\n$d = $word.ActiveDocument\n\n# 1) Swap the introductory wording in one shot with Find/Replace.\n$find = $d.Content.Find\n$find.Execute(\"Embargoed but can be requested from\", $false, $false, $false, $false, $false, $true, 1, $false, \"Available on the web at\", 2)\n\n# 2) Work out where the existing (repository) hyperlink currently starts - that is\n#    where we need to splice in the new link and the connecting words.\n$repoHyperlink = $d.Hyperlinks.Item(1)\n$repoStart = $repoHyperlink.Range.Start\n\n# 3) Insert the connecting words immediately before the repository hyperlink.\n$connectorText = \" or as a PDF from \"\n$connectorPoint = $d.Range($repoStart, $repoStart)\n$connectorPoint.InsertBefore($connectorText)\n\n# 4) Re-resolve the repository hyperlink's (now shifted) start position, then\n#    insert the new link's display/URL text immediately before the connecting\n#    words we just added.\n$repoHyperlink2 = $d.Hyperlinks.Item(1)\n$repoStart2 = $repoHyperlink2.Range.Start\n$newLinkText = \"https://benjarman.uk/phd_thesis\"\n$newLinkStart = $repoStart2 - $connectorText.Length\n$newLinkPoint = $d.Range($newLinkStart, $newLinkStart)\n$newLinkPoint.InsertBefore($newLinkText)\n\n# 5) Turn the freshly-inserted text into an actual hyperlink pointing at the new\n#    location.\n$newLinkSpan = $d.Range($newLinkStart, $newLinkStart + $newLinkText.Length)\n$d.Hyperlinks.Add($newLinkSpan, $newLinkText) | Out-Null\n"}
